# Updates the cryptocurrency price/volume table to reflect the latest
# scrape (GitHub Actions job "Updated cryptos list").
#
# Columns: A=index, B=Coin, C=Link, D=Price, E=Volume(1h)
# Price cells in column D that look like plain numbers ("309.79", "0.634", ...)
# are written with a leading apostrophe so Excel keeps them as literal text
# (matching how they are already stored in the workbook) instead of silently
# converting them to floating point numbers and losing the exact formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.951.86'
$ws.Range("E2").Value = '  +0.01%  '
$ws.Range("D3").Value = '2.294.47'
$ws.Range("E3").Value = '  +2.03%  '
$ws.Range("D4").Value = '''1.01'
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '''111.77'
$ws.Range("E5").Value = '  -3.69%  '
$ws.Range("D6").Value = '''309.79'
$ws.Range("E6").Value = '  +3.00%  '
$ws.Range("D7").Value = '''0.634'
$ws.Range("E7").Value = '  +1.07%  '
$ws.Range("E8").Value = '  -0.30%  '
$ws.Range("D9").Value = '''0.618'
$ws.Range("E9").Value = '  -1.07%  '
$ws.Range("D10").Value = '''44.39'
$ws.Range("E10").Value = '  -4.06%  '
$ws.Range("E11").Value = '  -0.77%  '
$ws.Range("D12").Value = '''8.87'
$ws.Range("E12").Value = '  -3.13%  '
$ws.Range("D13").Value = '''1.06'
$ws.Range("E13").Value = '  +18.58%  '
$ws.Range("E14").Value = '  -0.67%  '
$ws.Range("D15").Value = '''15.64'
$ws.Range("E15").Value = '  +1.65%  '
$ws.Range("D16").Value = '2.635.17'
$ws.Range("E16").Value = '  +1.77%  '
$ws.Range("D17").Value = '2.303.25'
$ws.Range("E17").Value = '  +1.74%  '
$ws.Range("D18").Value = '42.978.33'
$ws.Range("E18").Value = '  +0.31%  '
$ws.Range("E19").Value = '  -0.71%  '
$ws.Range("E20").Value = '  -5.72%  '
$ws.Range("D21").Value = '''76.03'
$ws.Range("E21").Value = '  +2.82%  '
$ws.Range("D22").Value = '''3.51'
$ws.Range("E22").Value = '  -6.15%  '
$ws.Range("D23").Value = '''2.50'
$ws.Range("E23").Value = '  +5.51%  '
$ws.Range("D24").Value = '''259.70'
$ws.Range("E24").Value = '  +11.60%  '
$ws.Range("D25").Value = '''9.03'
$ws.Range("E25").Value = '  -5.78%  '
$ws.Range("D26").Value = '''11.83'
$ws.Range("E26").Value = '  -3.06%  '
$ws.Range("E27").Value = '  -0.17%  '
$ws.Range("D28").Value = '''39.46'
$ws.Range("E28").Value = '  -2.12%  '
$ws.Range("E29").Value = '  +0.21%  '
$ws.Range("D30").Value = '''22.30'
$ws.Range("E30").Value = '  +4.46%  '
$ws.Range("D31").Value = '''173.40'
$ws.Range("E31").Value = '  -1.44%  '
$ws.Range("E32").Value = '  -2.90%  '
$ws.Range("E33").Value = '  -0.94%  '
$ws.Range("D34").Value = '''5.74'
$ws.Range("E34").Value = '  +0.51%  '
$ws.Range("D35").Value = '''5.14'
$ws.Range("E35").Value = '  +6.01%  '
$ws.Range("E36").Value = '  +1.12%  '
$ws.Range("D37").Value = '''4.19'
$ws.Range("E37").Value = '  -8.28%  '
$ws.Range("D38").Value = '''0.0378'
$ws.Range("E38").Value = '  +0.97%  '
$ws.Range("E39").Value = '  -2.18%  '
$ws.Range("D40").Value = '''2.62'
$ws.Range("E40").Value = '  +0.59%  '
$ws.Range("D41").Value = '''72.07'
$ws.Range("E41").Value = '  -0.61%  '
$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").Value = '''1.49'
$ws.Range("E42").Value = '  +10.57%  '
$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").Value = '''0.232'
$ws.Range("E43").Value = '  -3.78%  '
$ws.Range("E44").Value = '  -0.14%  '
$ws.Range("D45").Value = '''12.46'
$ws.Range("E45").Value = '  -7.98%  '
$ws.Range("D46").Value = '''5.72'
$ws.Range("E46").Value = '  +2.08%  '
$ws.Range("D47").Value = '''108.56'
$ws.Range("E47").Value = '  -0.26%  '
$ws.Range("D48").Value = '''8.96'
$ws.Range("E48").Value = '  +3.39%  '
$ws.Range("D49").Value = '''1.29'
$ws.Range("E49").Value = '  -4.95%  '
$ws.Range("D50").Value = '''0.0987'
$ws.Range("E50").Value = '  -0.48%  '
$ws.Range("B51").Value = 'MinaProtocolToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina'
$ws.Range("D51").Value = '''1.45'
$ws.Range("E51").Value = '  +20.87%  '
